$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "69.103.68"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "3.769.43"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "624.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("D7").Value = "3.767.79"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.460"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.72"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "4.408.54"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "3.754.66"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "69.093.71"
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.59"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.708"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("E24").Value = "  +4.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.21"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "3.920.31"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +3.51%  "
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.22"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.75"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.721.65"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.98"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.164"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +13.86%  "
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.44"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.02%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.966"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.297"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.65"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.65"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("E51").Value = "  -0.18%  "
